$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new values
$data = @(
    @(1, 3, 2, 7, 4, 4, 2, 23, 5),
    @(2, 2, 2, 7, 3, 5, 1, 12, 5),
    @(3, 1, 3, 2, 8, 1, 5, 56, 5),
    @(4, 4, 4, 6, 8, 2, 4, 45, 5),
    @(5, 2, 0, 5, 3, 3, 3, 34, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select() | Out-Null
